$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowsData = @(
    @("ECs", "Efna4", "Epha2", "ECs", 3, 1, 0.7471093333333333, 2.241328, 0.4648531405557055, 0.4648531405557054, 3, 1, 12.64154566666667, 37.924637, 0.5019409384896431, 0.5019409384896429, 9.444616755326223, 85.00155079793601, 0.2333288216303888, 0.2333288216303886),
    @("ECs", "Efna4", "Epha2", "FAPs", 3, 1, 0.7471093333333333, 2.241328, 0.4648531405557055, 0.4648531405557054, 3, 1, 1.291159666666667, 3.873479, 0.05126634921989955, 0.05126634921989955, 0.9646374377902222, 8.681736940112, 0.02383132343969585, 0.02383132343969584),
    @("ECs", "Efna4", "Epha2", "M1", 3, 1, 0.7471093333333333, 2.241328, 0.4648531405557055, 0.4648531405557054, 3, 1, 0.97959, 2.93877, 0.03889526936817373, 0.03889526936817372, 0.7318608318399999, 6.586747486559999, 0.01808058811855569, 0.01808058811855568),
    @("ECs", "Efna4", "Epha2", "M2", 3, 1, 0.7471093333333333, 2.241328, 0.4648531405557055, 0.4648531405557054, 3, 1, 1.079938333333333, 3.239815, 0.04287966636655804, 0.04287966636655803, 0.8068320082577777, 7.261488074319998, 0.01993274757647536, 0.01993274757647535),
    @("ECs", "Efna4", "Epha2", "Neutro", 3, 1, 0.7471093333333333, 2.241328, 0.4648531405557055, 0.4648531405557054, 2, 0.6666666666666666, 0.265959, 0.7978769999999999, 0.01056007814074267, 0.01056007814074267, 0.198700451184, 1.788304060656, 0.004908885488237888, 0.004908885488237886),
    @("ECs", "Efna4", "Epha2", "sCs", 3, 1, 0.7471093333333333, 2.241328, 0.4648531405557055, 0.4648531405557054, 3, 1, 8.927132333333335, 26.781397, 0.3544576984149831, 0.3544576984149831, 6.669543886135112, 60.025894975216, 0.164770774302352, 0.164770774302352),
    @("FAPs", "Efna4", "Epha2", "ECs", 3, 1, 0.4799286666666667, 1.439786, 0.2986127170267525, 0.2986127170267524, 3, 1, 12.64154566666667, 37.924637, 0.5019409384896431, 0.5019409384896429, 6.067040156409112, 54.60336140768201, 0.1498859474293504, 0.1498859474293503),
    @("FAPs", "Efna4", "Epha2", "FAPs", 3, 1, 0.4799286666666667, 1.439786, 0.2986127170267525, 0.2986127170267524, 3, 1, 1.291159666666667, 3.873479, 0.05126634921989955, 0.05126634921989955, 0.6196645372771111, 5.576980835494, 0.01530878383259654, 0.01530878383259653),
    @("FAPs", "Efna4", "Epha2", "M1", 3, 1, 0.4799286666666667, 1.439786, 0.2986127170267525, 0.2986127170267524, 3, 1, 0.97959, 2.93877, 0.03889526936817373, 0.03889526936817372, 0.47013332258, 4.23119990322, 0.01161462206551778, 0.01161462206551777),
    @("FAPs", "Efna4", "Epha2", "M2", 3, 1, 0.4799286666666667, 1.439786, 0.2986127170267525, 0.2986127170267524, 3, 1, 1.079938333333333, 3.239815, 0.04287966636655804, 0.04287966636655803, 0.5182933643988888, 4.664640279589999, 0.01280441367891855, 0.01280441367891855),
    @("FAPs", "Efna4", "Epha2", "Neutro", 3, 1, 0.4799286666666667, 1.439786, 0.2986127170267525, 0.2986127170267524, 2, 0.6666666666666666, 0.265959, 0.7978769999999999, 0.01056007814074267, 0.01056007814074267, 0.127641348258, 1.148772134322, 0.003153373625621987, 0.003153373625621986),
    @("FAPs", "Efna4", "Epha2", "sCs", 3, 1, 0.4799286666666667, 1.439786, 0.2986127170267525, 0.2986127170267524, 3, 1, 8.927132333333335, 26.781397, 0.3544576984149831, 0.3544576984149831, 4.284386717893556, 38.559480461042, 0.1058455763947473, 0.1058455763947473),
    @("sCs", "Efna4", "Epha2", "ECs", 3, 1, 0.3801563333333333, 1.140469, 0.2365341424175421, 0.2365341424175421, 3, 1, 12.64154566666667, 37.924637, 0.5019409384896431, 0.5019409384896429, 4.805763648305889, 43.251872834753, 0.118726169429904, 0.1187261694299039),
    @("sCs", "Efna4", "Epha2", "FAPs", 3, 1, 0.3801563333333333, 1.140469, 0.2365341424175421, 0.2365341424175421, 3, 1, 1.291159666666667, 3.873479, 0.05126634921989955, 0.05126634921989955, 0.4908425246278889, 4.417582721651, 0.01212624194760717, 0.01212624194760717),
    @("sCs", "Efna4", "Epha2", "M1", 3, 1, 0.3801563333333333, 1.140469, 0.2365341424175421, 0.2365341424175421, 3, 1, 0.97959, 2.93877, 0.03889526936817373, 0.03889526936817372, 0.37239734257, 3.35157608313, 0.00920005918410027, 0.009200059184100266),
    @("sCs", "Efna4", "Epha2", "M2", 3, 1, 0.3801563333333333, 1.140469, 0.2365341424175421, 0.2365341424175421, 3, 1, 1.079938333333333, 3.239815, 0.04287966636655804, 0.04287966636655803, 0.410545397026111, 3.694908573234999, 0.01014250511116413, 0.01014250511116413),
    @("sCs", "Efna4", "Epha2", "Neutro", 3, 1, 0.3801563333333333, 1.140469, 0.2365341424175421, 0.2365341424175421, 2, 0.6666666666666666, 0.265959, 0.7978769999999999, 0.01056007814074267, 0.01056007814074267, 0.101105998257, 0.9099539843129999, 0.002497819026882802, 0.0024978190268828),
    @("sCs", "Efna4", "Epha2", "sCs", 3, 1, 0.3801563333333333, 1.140469, 0.2365341424175421, 0.2365341424175421, 3, 1, 8.927132333333335, 26.781397, 0.3544576984149831, 0.3544576984149831, 3.393705895021445, 30.543353055193, 0.08384134771788382, 0.08384134771788379)
)

$r = 2
foreach ($row in $rowsData) {
    $c = 1
    foreach ($val in $row) {
        $ws.Cells.Item($r, $c).Value = $val
        $c++
    }
    $r++
}

